$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-11-18 01:20:10"

# Insert a new row at position 9, pushing existing rows 9-11 down to rows 10-12
$ws.Rows.Item(9).EntireRow.Insert()

# Refresh the "取得日時" (fetched-at) timestamp for the rows that already existed
# and kept their row numbers (rows 2-8).
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# Fill the newly inserted row 9 with the new job posting that was appended in this run.
$ws.Range("A9").Value = $newTimestamp
$ws.Range("B9").Value = "【技術パートナー募集】リード獲得・育成システム構築"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5436021"
$ws.Range("G9").Value = 33

# Rows 10, 11 and 12 now hold the data that used to live in rows 9, 10 and 11;
# only their timestamp needs to be refreshed to the new run time.
$ws.Range("A10").Value = $newTimestamp
$ws.Range("A11").Value = $newTimestamp
$ws.Range("A12").Value = $newTimestamp

# The engine does not shift hyperlink ranges/targets when rows are inserted, so the
# hyperlink-to-target mapping for column F is now stale for rows 9-11 and missing for
# the new row 12. Rebuild the whole hyperlinks collection from scratch so every F cell
# points at the correct URL (this also re-applies the blue/underlined Hyperlink style).
$ws.Hyperlinks.Delete()

$links = @(
    @{ Cell = "F2"; Url = "https://www.lancers.jp/work/detail/5428507" },
    @{ Cell = "F3"; Url = "https://www.lancers.jp/work/detail/5423720" },
    @{ Cell = "F4"; Url = "https://www.lancers.jp/work/detail/5434977" },
    @{ Cell = "F5"; Url = "https://www.lancers.jp/work/detail/5416328" },
    @{ Cell = "F6"; Url = "https://www.lancers.jp/work/detail/5419380" },
    @{ Cell = "F7"; Url = "https://www.lancers.jp/work/detail/5435875" },
    @{ Cell = "F8"; Url = "https://www.lancers.jp/work/detail/5431107" },
    @{ Cell = "F9"; Url = "https://www.lancers.jp/work/detail/5436021" },
    @{ Cell = "F10"; Url = "https://www.lancers.jp/work/detail/5429882" },
    @{ Cell = "F11"; Url = "https://www.lancers.jp/work/detail/5435667" },
    @{ Cell = "F12"; Url = "https://www.lancers.jp/work/detail/5435519" }
)

foreach ($link in $links) {
    $ws.Hyperlinks.Add($ws.Range($link.Cell), $link.Url)
    $ws.Range($link.Cell).Style = "Hyperlink"
}
